# Auto-generated Excel COM-interop script to update cryptos worksheet
# Commit: Updated cryptos list on Tue Apr 25 05:59:16 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.489.76"
Set-TextValue $ws.Range("E2") "  -1.51%  "
Set-TextValue $ws.Range("D3") "1.834.58"
Set-TextValue $ws.Range("E3") "  -2.17%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  -0.77%  "
Set-TextValue $ws.Range("D5") "331.43"
Set-TextValue $ws.Range("E5") "  -1.20%  "
Set-TextValue $ws.Range("E6") "  -0.85%  "
Set-TextValue $ws.Range("D7") "0.4619"
Set-TextValue $ws.Range("E7") "  -2.89%  "
Set-TextValue $ws.Range("D8") "0.3832"
Set-TextValue $ws.Range("E8") "  -2.61%  "
Set-TextValue $ws.Range("E9") "  -0.25%  "
Set-TextValue $ws.Range("D10") "0.07924"
Set-TextValue $ws.Range("E10") "  -0.89%  "
Set-TextValue $ws.Range("D11") "0.9737"
Set-TextValue $ws.Range("E11") "  -3.99%  "
Set-TextValue $ws.Range("D12") "21.11"
Set-TextValue $ws.Range("E12") "  -3.18%  "
Set-TextValue $ws.Range("D13") "1.841.74"
Set-TextValue $ws.Range("E13") "  -2.69%  "
Set-TextValue $ws.Range("D14") "5.901"
Set-TextValue $ws.Range("E14") "  -2.10%  "
Set-TextValue $ws.Range("D15") "7.061"
Set-TextValue $ws.Range("E15") "  -1.75%  "
Set-TextValue $ws.Range("D16") "1.001"
Set-TextValue $ws.Range("E16") "  -1.30%  "
Set-TextValue $ws.Range("D17") "88.11"
Set-TextValue $ws.Range("E17") "  -0.17%  "
Set-TextValue $ws.Range("D18") "0.06644"
Set-TextValue $ws.Range("E18") "  -1.02%  "
Set-TextValue $ws.Range("D19") "0.00001029"
Set-TextValue $ws.Range("E19") "  -1.84%  "
Set-TextValue $ws.Range("D20") "17.06"
Set-TextValue $ws.Range("E20") "  +0.16%  "
Set-TextValue $ws.Range("E21") "  -0.81%  "
Set-TextValue $ws.Range("D22") "27.475.92"
Set-TextValue $ws.Range("E22") "  -1.55%  "
Set-TextValue $ws.Range("D23") "5.356"
Set-TextValue $ws.Range("E23") "  -2.53%  "
Set-TextValue $ws.Range("D24") "10.83"
Set-TextValue $ws.Range("E24") "  -1.17%  "
Set-TextValue $ws.Range("D25") "2.308"
Set-TextValue $ws.Range("E25") "  -1.31%  "
Set-TextValue $ws.Range("D26") "2.045.11"
Set-TextValue $ws.Range("E26") "  -3.15%  "
Set-TextValue $ws.Range("D27") "157.40"
Set-TextValue $ws.Range("E27") "  -0.60%  "
Set-TextValue $ws.Range("D28") "19.42"
Set-TextValue $ws.Range("E28") "  -1.97%  "
Set-TextValue $ws.Range("D29") "2.072"
Set-TextValue $ws.Range("E29") "  -1.06%  "
Set-TextValue $ws.Range("D30") "5.311"
Set-TextValue $ws.Range("E30") "  -2.66%  "
Set-TextValue $ws.Range("D31") "119.26"
Set-TextValue $ws.Range("E31") "  -1.61%  "
Set-TextValue $ws.Range("D32") "0.9573"
Set-TextValue $ws.Range("E32") "  -1.60%  "
Set-TextValue $ws.Range("D33") "0.09306"
Set-TextValue $ws.Range("E33") "  -2.26%  "
Set-TextValue $ws.Range("D34") "3.581"
Set-TextValue $ws.Range("E34") "  -1.36%  "
Set-TextValue $ws.Range("D35") "5.255"
Set-TextValue $ws.Range("E35") "  -1.32%  "
Set-TextValue $ws.Range("D36") "1.317"
Set-TextValue $ws.Range("E36") "  -2.34%  "
Set-TextValue $ws.Range("D37") "0.05937"
Set-TextValue $ws.Range("E37") "  -2.32%  "
Set-TextValue $ws.Range("D38") "0.02196"
Set-TextValue $ws.Range("E38") "  -1.49%  "
Set-TextValue $ws.Range("D39") "1.165"
Set-TextValue $ws.Range("E39") "  -3.34%  "
Set-TextValue $ws.Range("D40") "8.071"
Set-TextValue $ws.Range("E40") "  -1.14%  "
Set-TextValue $ws.Range("D41") "0.5802"
Set-TextValue $ws.Range("E41") "  -2.73%  "
Set-TextValue $ws.Range("D42") "0.1841"
Set-TextValue $ws.Range("E42") "  -2.69%  "
Set-TextValue $ws.Range("D43") "10.07"
Set-TextValue $ws.Range("E43") "  -2.39%  "
Set-TextValue $ws.Range("D44") "1.272"
Set-TextValue $ws.Range("E44") "  +1.02%  "
Set-TextValue $ws.Range("B45") "EnergySwap"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "12.05"
Set-TextValue $ws.Range("E45") "  -0.59%  "
Set-TextValue $ws.Range("B46") "Decentraland"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.5492"
Set-TextValue $ws.Range("E46") "  -2.99%  "
Set-TextValue $ws.Range("D47") "1.874"
Set-TextValue $ws.Range("E47") "  -2.85%  "
Set-TextValue $ws.Range("D48") "0.06652"
Set-TextValue $ws.Range("E48") "  -1.97%  "
Set-TextValue $ws.Range("D49") "110.54"
Set-TextValue $ws.Range("E49") "  -1.44%  "
Set-TextValue $ws.Range("D50") "1.042"
Set-TextValue $ws.Range("E50") "  -2.46%  "
Set-TextValue $ws.Range("D51") "1.001"
Set-TextValue $ws.Range("E51") "  -0.96%  "
